# Updated cryptos list on Sun May 26 03:25:16 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# swaps the two row-pairs whose rank order changed (dogwifhat <-> FirstDigitalUSD,
# Bittensor <-> OKB). All D/E cells are plain text (not numbers), so for any
# new value that Excel would otherwise auto-parse as a number (losing a
# trailing zero, e.g. "1.00" -> 1, or switching to scientific notation) we
# briefly force the cell to Text format, assign the literal string, then put
# the cell's style back to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.953.92'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '3.738.98'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.73%  '
$ws.Range("D7").Value = '3.735.83'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  +4.18%  '
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.28%  '
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").Value = '4.365.91'
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '3.734.96'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '69.009.70'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.42'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.45%  '
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '490.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.79%  '
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000147'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.61%  '
$ws.Range("E26").Value = '  -1.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.98%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.43'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.25%  '
$ws.Range("D34").Value = '3.883.88'
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("D35").Value = '3.674.80'
$ws.Range("E35").Value = '  +0.27%  '
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.137'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.03%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.60%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '426.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.91%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '39.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.15%  '
$ws.Range("D51").Value = '2.781.33'
$ws.Range("E51").Value = '  -0.17%  '
